$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("L2").Value = "*maa://24633 (55.48), *maa://30515 (69.0), *maa://34787 (71.21), ***maa://20792 (11.93), maa://39402 (86.49), ***maa://29083 (27.78)"
$ws.Range("X3").Value = "maa://27396 (84.69), maa://27484 (96.0), maa://27480 (82.86)"
$ws.Range("AF4").Value = "*maa://30062 (60.0), ***maa://26209 (13.04), *maa://39394 (68.18)"
$ws.Range("T7").Value = "maa://21291 (85.37)"
$ws.Range("A8").Value = "更新日期：2024.12.07 23:48:34"
$ws.Range("D10").Value = "***maa://25695 (19.1), **maa://32237 (41.46), ***maa://34206 (18.18), ***maa://39951 (17.5), ***maa://39243 (28.57)"
$ws.Range("D11").Value = "maa://36707 (99.4)"
$ws.Range("T11").Value = "maa://22747 (93.33), maa://22501 (98.51)"
$ws.Range("AF13").Value = "**maa://22737 (31.88), maa://39883 (90.57), *maa://39885 (56.0)"
$ws.Range("D15").Value = "*maa://22743 (76.84), maa://22734 (83.9), *maa://30808 (63.93), ***maa://36048 (28.57)"
$ws.Range("H17").Value = "maa://22430 (88.59), maa://39599 (85.71)"
$ws.Range("D18").Value = "maa://24570 (97.04)"
$ws.Range("L18").Value = "maa://22466 (88.65), *maa://22732 (50.6)"
$ws.Range("AB19").Value = "*maa://30709 (62.94), *maa://36668 (54.67)"
$ws.Range("P23").Value = "maa://30587 (91.85), *maa://29748 (75.59), ***maa://29785 (16.42), *maa://37566 (71.43)"
$ws.Range("P26").Value = "*maa://30968 (63.16), maa://39870 (100.0)"
$ws.Range("X28").Value = "maa://39929 (89.4), ***maa://39723 (14.29), maa://41749 (90.0)"
$ws.Range("AB30").Value = "maa://42979 (97.87)"
$ws.Range("L31").Value = "maa://35926 (93.89), maa://36258 (82.42), maa://43904 (100.0)"
$ws.Range("H32").Value = "maa://21895 (97.22), maa://36667 (98.31), **maa://20793 (38.78), maa://22760 (100.0)"
$ws.Range("T32").Value = "maa://42859 (95.24), maa://41108 (87.76), maa://41238 (96.05)"
$ws.Range("L35").Value = "maa://41296 (96.3)"
$ws.Range("P39").Value = "maa://24709 (91.87)"
$ws.Range("H41").Value = "maa://24466 (93.33)"
$ws.Range("H46").Value = "maa://35931 (92.42), *maa://43901 (80.0)"
$ws.Range("H47").Value = "maa://27410 (96.15), maa://29661 (97.78), maa://28038 (84.62)"
$ws.Range("H53").Value = "maa://32534 (93.49), **maa://32434 (34.78)"
$ws.Range("H55").Value = "maa://32532 (92.58)"
